$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - use a leading apostrophe via Formula so that
# numeric-looking text values (e.g. "579.12") are stored as text, matching
# the original inline-string cell type rather than being parsed as numbers.
$priceUpdates = [ordered]@{
    "D2" = "62.627.38"
    "D3" = "2.566.87"
    "D5" = "579.12"
    "D6" = "143.45"
    "D8" = "0.587"
    "D10" = "5.53"
    "D13" = "26.80"
    "D14" = "3.027.51"
    "D15" = "62.568.85"
    "D17" = "2.566.48"
    "D18" = "11.07"
    "D19" = "337.98"
    "D23" = "66.83"
    "D27" = "0.998"
    "D28" = "7.90"
    "D29" = "8.19"
    "D30" = "1.92"
    "D31" = "453.76"
    "D32" = "0.0₃0793"
    "D33" = "176.37"
    "D34" = "1.62"
    "D37" = "18.75"
    "D38" = "4.43"
    "D41" = "40.28"
    "D42" = "157.31"
    "D45" = "20.87"
    "D49" = "17.93"
}

foreach ($cell in $priceUpdates.Keys) {
    $ws.Range($cell).Formula = "'" + $priceUpdates[$cell]
}

# Volume(1h) column (E) updates - plain text assignment (the percent sign and
# surrounding spaces already prevent Excel from treating these as numbers).
$volumeUpdates = [ordered]@{
    "E2" = "  -0.85%  "
    "E3" = "  +0.35%  "
    "E4" = "  -0.07%  "
    "E5" = "  -0.55%  "
    "E6" = "  -3.00%  "
    "E7" = "  -0.05%  "
    "E8" = "  +0.23%  "
    "E9" = "  -2.23%  "
    "E10" = "  -1.07%  "
    "E11" = "  -0.50%  "
    "E12" = "  -1.89%  "
    "E13" = "  -3.22%  "
    "E14" = "  +0.24%  "
    "E15" = "  -0.78%  "
    "E16" = "  -2.06%  "
    "E17" = "  +0.07%  "
    "E18" = "  -2.92%  "
    "E19" = "  -0.62%  "
    "E20" = "  -1.91%  "
    "E21" = "  -2.69%  "
    "E22" = "  -0.01%  "
    "E23" = "  +0.93%  "
    "E24" = "  -4.22%  "
    "E25" = "  +1.31%  "
    "E26" = "  -4.31%  "
    "E27" = "  -0.39%  "
    "E28" = "  -2.86%  "
    "E29" = "  -4.00%  "
    "E30" = "  -1.62%  "
    "E31" = "  +3.39%  "
    "E32" = "  -3.86%  "
    "E33" = "  -0.47%  "
    "E34" = "  +1.05%  "
    "E35" = "  -0.03%  "
    "E36" = "  -3.02%  "
    "E37" = "  -2.75%  "
    "E38" = "  -2.17%  "
    "E39" = "  +0.00%  "
    "E40" = "  -3.77%  "
    "E41" = "  +1.19%  "
    "E42" = "  +3.86%  "
    "E43" = "  -4.09%  "
    "E44" = "  +3.05%  "
    "E45" = "  -1.47%  "
    "E46" = "  -3.14%  "
    "E47" = "  -2.17%  "
    "E48" = "  -3.34%  "
    "E49" = "  -2.68%  "
    "E50" = "  +0.28%  "
    "E51" = "  -4.32%  "
}

foreach ($cell in $volumeUpdates.Keys) {
    $ws.Range($cell).Value2 = $volumeUpdates[$cell]
}

